$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (the current most-recent week). This shifts
# the existing data rows 2-68 down to 3-69, preserving all of their values.
$ws.Rows.Item(2).Insert()

# The freshly inserted row 2 doesn't carry the date-column number format
# that the rest of the table uses, so copy formats down from row 3 (which
# now holds what used to be row 2) onto the new row 2.
$ws.Range("A3:AA3").Copy()
$ws.Range("A2:AA2").PasteSpecial(-4122)

# Fill in the new week's data in row 2.
$ws.Range("A2").Value = 43920
$ws.Range("B2").Value = 43924
$ws.Range("C2").Value = 255.7
$ws.Range("D2").Value = 263.33
$ws.Range("E2").Value = 247.6
$ws.Range("F2").Value = 248.89
$ws.Range("G2").Value = 274.065
$ws.Range("H2").Value = 232.77500000000001

# Match the saved selection.
$ws.Range("F2").Select()
